$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the sheet with the initial data
$ws.Range("A1").Value = "lkdfmkf"
$ws.Range("A2").Value = "njf"
$ws.Range("B2").Value = "nkfjnn"
$ws.Range("C3").Value = "jnnjvf"
$ws.Range("D4").Value = "jnefvn"

# Leave the active selection on D4, matching the last-edited cell
$ws.Range("D4").Select() | Out-Null
